$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.034.61'
$ws.Range("E2").Value = '  -0.15%  '
$ws.Range("D3").Value = '3.422.51'
$ws.Range("E3").Value = '  +0.10%  '
$ws.Range("E4").Value = '  +0.13%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '409.79'
$ws.Range("E5").Value = '  +0.98%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '129.04'
$ws.Range("E6").Value = '  -2.60%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.628'
$ws.Range("E7").Value = '  +6.65%  '
$ws.Range("E8").Value = '  -0.06%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.749'
$ws.Range("E9").Value = '  +11.18%  '
$ws.Range("E10").Value = '  +19.78%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '42.60'
$ws.Range("E11").Value = '  +0.29%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0000226'
$ws.Range("E12").Value = '  +75.14%  '
$ws.Range("E13").Value = '  -0.33%  '
$ws.Range("D14").Value = '3.963.74'
$ws.Range("E14").Value = '  +1.11%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '8.98'
$ws.Range("E15").Value = '  +6.63%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '21.05'
$ws.Range("E16").Value = '  +6.02%  '
$ws.Range("D17").Value = '3.430.34'
$ws.Range("E17").Value = '  +0.94%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '12.45'
$ws.Range("E18").Value = '  +13.39%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.07'
$ws.Range("D20").Value = '61.998.88'
$ws.Range("E20").Value = '  +0.06%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '399.17'
$ws.Range("E21").Value = '  +26.04%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '89.44'
$ws.Range("E22").Value = '  +6.01%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.19'
$ws.Range("E23").Value = '  -0.82%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '13.22'
$ws.Range("E24").Value = '  +3.40%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.24'
$ws.Range("E25").Value = '  +4.03%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '32.69'
$ws.Range("E26").Value = '  +10.24%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.74'
$ws.Range("E27").Value = '  +6.00%  '
$ws.Range("E28").Value = '  +0.24%  '
$ws.Range("B29").Value = 'Toncoin'
$ws.Range("C29").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.75'
$ws.Range("E29").Value = '  +2.02%  '
$ws.Range("B30").Value = 'RenderToken'
$ws.Range("C30").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.58'
$ws.Range("E30").Value = '  -2.03%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.119'
$ws.Range("E31").Value = '  +2.60%  '
$ws.Range("E32").Value = '  -0.21%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '11.86'
$ws.Range("E33").Value = '  +4.24%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '43.07'
$ws.Range("E34").Value = '  +2.46%  '
$ws.Range("E35").Value = '  +0.64%  '
$ws.Range("E36").Value = '  +3.62%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '53.80'
$ws.Range("E37").Value = '  +3.69%  '
$ws.Range("E38").Value = '  +0.01%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.37'
$ws.Range("E39").Value = '  -1.75%  '
$ws.Range("E40").Value = '  +7.28%  '
$ws.Range("E41").Value = '  -3.41%  '
$ws.Range("E42").Value = '  +7.03%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '141.73'
$ws.Range("E43").Value = '  +1.98%  '
$ws.Range("E44").Value = '  -1.02%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '4.13'
$ws.Range("E45").Value = '  +4.14%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.42'
$ws.Range("E46").Value = '  +8.98%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '16.69'
$ws.Range("E47").Value = '  -0.20%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '21.99'
$ws.Range("E48").Value = '  +3.30%  '
$ws.Range("D49").Value = '2.121.30'
$ws.Range("E49").Value = '  -0.36%  '
$ws.Range("B50").Value = 'ThetaToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.03'
$ws.Range("E50").Value = '  +8.20%  '
$ws.Range("B51").Value = 'ApeXProtocol'
$ws.Range("C51").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.37'
$ws.Range("E51").Value = '  +2.60%  '
